# ADD: Lab 2 Sweep Data and code
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header for column D: "Calculated Distance"
$ws.Range("D1").Value = "Calculated Distance"

# Rows 2 & 3 did not have usable calculated distances
$ws.Range("D2").Value = "Did not include"
$ws.Range("D3").Value = "Did not include"

# Rows 4-11: back-calculate distance from the voltage reading using the
# inverse of the exponential calibration curve fit (D14/D15 below).
$ws.Range("D4").Formula = "=(LN(C4) - LN(4.2223))/-0.024"
$ws.Range("D5").Formula = "=(LN(C5) - LN(4.2223))/-0.024"
$ws.Range("D6").Formula = "=(LN(C6) - LN(4.2223))/-0.024"
$ws.Range("D7").Formula = "=(LN(C7) - LN(4.2223))/-0.024"
$ws.Range("D8").Formula = "=(LN(C8) - LN(4.2223))/-0.024"
$ws.Range("D9").Formula = "=(LN(C9) - LN(4.2223))/-0.024"
$ws.Range("D10").Formula = "=(LN(C10) - LN(4.2223))/-0.024"
$ws.Range("D11").Formula = "=(LN(C11) - LN(4.2223))/-0.024"

# Move "Chart 1" (the sweep-data scatter chart) over to make room for the
# new "Calculated Distance" column/data.
$chartObj = $ws.ChartObjects(1)
$chartObj.Left = 655.3662109375
$chartObj.Top = 19.12496062992126
$chartObj.Width = 433.0625
$chartObj.Height = 216

# Leave selection on E3, matching where editing left off.
$ws.Range("E3").Select()
